$wb = $excel.ActiveWorkbook

# --- USERACCOUNTMANAGEMENTDATA sheet: add two data rows for the new test case ---
$ws3 = $wb.Worksheets.Item("USERACCOUNTMANAGEMENTDATA")

$ws3.Range("A10").Value = "verifyThatUserCannotChangePasswordWhenAllFieldsAreEmpty"
$ws3.Range("B10").Value = "yes"
$ws3.Range("C10").Value = "'Admin"
$ws3.Range("D10").Value = "'admin123"
$ws3.Range("E10").Value = "'Sunil"
$ws3.Range("F10").Value = "'chrome"
$ws3.Range("G10").Value = "'"
$ws3.Range("H10").Value = "'"
$ws3.Range("I10").Value = "'"

$ws3.Range("A11").Value = "verifyThatUserCannotChangePasswordWhenAllFieldsAreEmpty"
$ws3.Range("B11").Value = "yes"
$ws3.Range("C11").Value = "Admin"
$ws3.Range("D11").Value = "admin123"
$ws3.Range("E11").Value = "Sunil"
$ws3.Range("F11").Value = "firefox"
$ws3.Range("G11").Value = "'"
$ws3.Range("H11").Value = "'"
$ws3.Range("I11").Value = "'"

$ws3.Range("A11").Select()

# --- RUNMANAGER sheet: register the new test case ---
$ws1 = $wb.Worksheets.Item("RUNMANAGER")

$ws1.Range("A10").Value = "verifyThatUserCannotChangePasswordWhenAllFieldsAreEmpty"
$ws1.Range("B10").Value = "To check this test is executed"
$ws1.Range("C10").Value = "yes"
$ws1.Range("D10").Value = "'9"
$ws1.Range("E10").Value = "'1"

# Leave RUNMANAGER as the active sheet/selection, as in the final workbook state
$ws1.Range("F10").Select()
